$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original data (A1:C6) before moving it
$values = $ws.Range("A1:C6").Value2

# Clear the old range
$ws.Range("A1:C6").ClearContents()

# Write the data shifted: right by 1 column (A->B), down by 3 rows (1->4)
$ws.Range("B4:D9").Value2 = $values

# New empty styled cell at H7 (underline font, same style as H12)
$ws.Range("H7").Font.Underline = $true

# New empty styled cell at N12 (same style as H12)
$ws.Range("N12").Font.Underline = $true

# Update selection to match target
$ws.Range("H7").Select()
